# Apply the "inserted mapping values #492" edit to the mapping workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "mapping"

# Update the id_start_year_usage (column C) values for the rows that changed.
$newValues = @{
    5  = 3
    6  = 3
    7  = 3
    8  = 3
    9  = 3
    10 = 3
    15 = 2
    16 = 2
    18 = 2
    19 = 2
    24 = 3
    26 = 2
    27 = 2
    33 = 2
    35 = 2
    36 = 2
    39 = 2
    40 = 2
    41 = 2
    42 = 2
    44 = 2
    46 = 2
    49 = 2
    51 = 2
    54 = 2
    55 = 2
    56 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

# Remove the last data row (row 65), which is no longer present after the edit.
$ws.Rows.Item(65).Delete()

# Match the saved selection state from the edit.
$ws.Range("C2:C64").Select()
